$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates (row 2) ---
$wsSchedule.Range("E2").Value = 781.7284702500003
$wsSchedule.Range("F2").Value = 12.92540460069445

# --- Detailed sheet updates ---
$wsDetailed.Range("B2").Value = 57.06003
$wsDetailed.Range("B3").Value = 57.06003
$wsDetailed.Range("B4").Value = 57.06003
$wsDetailed.Range("B5").Value = 57.09
$wsDetailed.Range("C5").Value = "historical"
$wsDetailed.Range("B6").Value = 56.97996
$wsDetailed.Range("C6").Value = "historical"
$wsDetailed.Range("B7").Value = 56.97996
$wsDetailed.Range("B8").Value = 56.97996
$wsDetailed.Range("B11").Value = 65
$wsDetailed.Range("B12").Value = 57.3
$wsDetailed.Range("B13").Value = 65
$wsDetailed.Range("B14").Value = 57.3
$wsDetailed.Range("B15").Value = 36.2
$wsDetailed.Range("B16").Value = 36.06029
$wsDetailed.Range("B17").Value = 18.4525
$wsDetailed.Range("B18").Value = 0.7
$wsDetailed.Range("B19").Value = 36.06011
$wsDetailed.Range("B20").Value = 36.06054
$wsDetailed.Range("B21").Value = 36.06011
$wsDetailed.Range("B22").Value = 46.43705
$wsDetailed.Range("B24").Value = 44.73383
$wsDetailed.Range("B27").Value = 36.06045
$wsDetailed.Range("B28").Value = 36.06045
$wsDetailed.Range("B29").Value = 36.06029
$wsDetailed.Range("B30").Value = 30.01081
$wsDetailed.Range("B35").Value = 25.55134
$wsDetailed.Range("B36").Value = -0.32102
$wsDetailed.Range("B37").Value = -3.01589
$wsDetailed.Range("B38").Value = -2.85215
$wsDetailed.Range("B39").Value = -2.79911
$wsDetailed.Range("B40").Value = 0.0109
$wsDetailed.Range("B41").Value = 9.38786
$wsDetailed.Range("B42").Value = 29.66317
$wsDetailed.Range("B43").Value = 9.42876
$wsDetailed.Range("B44").Value = 9.52481
$wsDetailed.Range("B45").Value = 9.754619999999999
$wsDetailed.Range("B46").Value = 36.06045
